$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.453.30"
$ws.Range("E2").Value = "  +1.05%  "
$ws.Range("D3").Value = "1.853.94"
$ws.Range("E3").Value = "  +1.35%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.03%  "
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4752"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.09%  "
$ws.Range("E8").Value = "  +2.79%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06349"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.02"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +12.76%  "
$ws.Range("D11").Value = "1.888.35"
$ws.Range("E11").Value = "  +3.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07472"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.20%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.974"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.69%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "84.79"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6246"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.26%  "
$ws.Range("D16").Value = "30.415.17"
$ws.Range("E16").Value = "  +1.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "245.11"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +8.47%  "
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.70"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.71%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007351"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.932"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.917"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.17%  "
$ws.Range("B24").Value = "Monero"
$ws.Range("C24").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "164.33"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.080"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.882"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1027"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.04%  "
$ws.Range("E29").Value = "  -1.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.053"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.844"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.04840"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.131"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.96%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.6997"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.27%  "
$ws.Range("E35").Value = "  +0.44%  "
$ws.Range("E36").Value = "  +5.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.683"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.8794"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.992"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.73%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "106.70"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.67%  "
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4075"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.505"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.91%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.191"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.70%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "63.66"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1203"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.38%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "33.99"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.620"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.71%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05502"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.25%  "
$ws.Range("E50").Value = "  -0.55%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3696"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.46%  "
